# Add a new "water" column (Q) to the sheet, with a header and values for
# rows 2-30 (row 31 has no data and is intentionally skipped).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("Q1").Value = "water"

# Values for column Q, rows 2-30
$values = @{
    2  = 4
    3  = 12
    4  = 7.2
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 6.5
    19 = 0
    20 = 7.3
    21 = 6.2
    22 = 4.5999999999999996
    23 = 6.6
    24 = 13
    25 = 11.5
    26 = 0
    27 = 10.5
    28 = 13
    29 = 0
    30 = 0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 17).Value = $values[$row]
}

# Update the active selection to match the edited area
$ws.Range("Q30").Select()
